$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 6 new rows above the "Regex" row (old row 19), pushing the
#     old rows 19 ("Regex") and 20 (blank spacer) down to 25 and 26.
$ws.Range("A19:C24").Insert(-4121)

# Excel copies formatting from the row above on insert, but the new rows
# in the real workbook keep the standard 30pt data-row height, so fix it.
$ws.Rows.Item(19).RowHeight = 30
$ws.Rows.Item(20).RowHeight = 30
$ws.Rows.Item(21).RowHeight = 30
$ws.Rows.Item(22).RowHeight = 30
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(24).RowHeight = 30

# --- Fill in the new DMS config rows (values set in the same order the
#     author typed them so new shared-string entries line up).
$ws.Range("A19").Value() = "DMSHomeUrl"
$ws.Range("A21").Value() = "DMSEmailTitle"
$ws.Range("A23").Value() = "DMSExcelReturnTitle"
$ws.Range("B23").Value() = "Compliance Waste Returns {1} {2}"
$ws.Range("B21").Value() = "Compliance Waste Returns {1} {2} - Email and Submission"
$ws.Range("A22").Value() = "DMSEmailDisclosureStatus"
$ws.Range("B22").Value() = "Internal Only"
$ws.Range("A24").Value() = "DMSExcelReturnDisclosureStatus"
$ws.Range("B24").Value() = "Public Register"
$ws.Range("C21").Value() = "Title to rename email file once uploaded to DMS"
$ws.Range("C23").Value() = "Title to rename excel return file once uploaded to DMS"
$ws.Range("C24").Value() = "Disclosure status for excel return file once uploaded to DMS"
$ws.Range("C22").Value() = "Disclosure status for email file once uploaded to DMS"
$ws.Range("C19").Value() = "Home url site for DMS"
$ws.Range("A20").Value() = "DMSPermitFolderUrl"
$ws.Range("C20").Value() = "Url for permit folder  on DMS"

# Rows 19 & 20 only have Name/Description (no Value) - Insert() copied a
# blank styled B-cell down from the row above, so drop it to match rows
# that genuinely only ever had two populated columns.
$ws.Range("B19").Clear()
$ws.Range("B20").Clear()

# --- Grow Table1 so it covers the new rows too.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C26"))

# --- Column A got a touch wider to fit the new, longer setting names
#     (e.g. "DMSExcelReturnDisclosureStatus").
$ws.Columns.Item(1).ColumnWidth = 31.25

# --- Restore the view: scrolled down with A21 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("A21").Select()
